# Auto-generated data refresh for Titan_Profits market-data columns (H:N)
# Updates per-sheet leve profit calculations (currentAveragePrice, LevePrice*, LeveProfit*)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H39" = 151.84616
    "I39" = 75
    "J39" = 324.75
    "K39" = 225
    "L39" = 974.25
    "M39" = 71
    "N39" = -1566.25
    "H133" = 13414.286
    "J133" = 13414.286
    "L133" = 13414.286
    "N133" = -23534.286
    "H138" = 3860325.5
    "I138" = 2179418
    "J138" = 4222039.5
    "K138" = 6538254
    "L138" = 12666118.5
    "M138" = -6533114
    "N138" = -12676398.5
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H32" = 15210.284
    "I32" = 1478.942
    "J32" = 204702.8
    "K32" = 1478.942
    "L32" = 204702.8
    "M32" = -1191.942
    "N32" = -205276.8
    "H41" = 21333.334
    "I41" = 2000
    "J41" = 60000
    "K41" = 2000
    "L41" = 60000
    "M41" = -1586
    "N41" = -60828
    "H45" = 580
    "I45" = 484.16666
    "J45" = 963.3333
    "K45" = 484.16666
    "L45" = 963.3333
    "M45" = -107.16666
    "N45" = -1717.3333
    "H60" = 50000
    "I60" = 50000
    "K60" = 50000
    "M60" = -49267
    "H61" = 2468.1333
    "I61" = 1390.1
    "J61" = 4624.2
    "K61" = 1390.1
    "L61" = 4624.2
    "M61" = -1178.1
    "N61" = -5048.2
    "H74" = 7263.811
    "I74" = 1251.8889
    "J74" = 23496
    "K74" = 1251.8889
    "L74" = 23496
    "M74" = -377.8888999999999
    "N74" = -25244
    "H77" = 7263.811
    "I77" = 1251.8889
    "J77" = 23496
    "K77" = 6259.4445
    "L77" = 117480
    "M77" = -1891.4445
    "N77" = -126216
    "H110" = 511.08
    "I110" = 485.31818
    "K110" = 485.31818
    "M110" = 1559.68182
    "H136" = 2468.1333
    "I136" = 1390.1
    "J136" = 4624.2
    "K136" = 4170.299999999999
    "L136" = 13872.6
    "M136" = -1620.299999999999
    "N136" = -18972.6
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H59" = 42800
    "J59" = 47360
    "L59" = 47360
    "N59" = -49054
    "H107" = 997.17645
    "I107" = 667.1
    "J107" = 1468.7142
    "K107" = 667.1
    "L107" = 1468.7142
    "M107" = 1252.9
    "N107" = -5308.7142
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H31" = 1030.8485
    "I31" = 839.29034
    "K31" = 839.29034
    "M31" = -544.29034
    "H34" = 1030.8485
    "I34" = 839.29034
    "K34" = 839.29034
    "M34" = -637.29034
    "H132" = 3295.318
    "I132" = 2734.625
    "J132" = 4790.5
    "K132" = 8203.875
    "L132" = 14371.5
    "M132" = -5673.875
    "N132" = -19431.5
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$values = @{
    "H12" = 38.333332
    "I12" = 14.5
    "J12" = 86
    "K12" = 43.5
    "L12" = 258
    "M12" = 129.5
    "N12" = -604
    "H76" = 3499.75
    "J76" = 3666.6667
    "L76" = 11000.0001
    "H79" = 3499.75
    "J79" = 3666.6667
    "L79" = 11000.0001
    "H80" = 1200
    "I80" = 0
    "J80" = 1200
    "K80" = 0
    "L80" = 3600
    "N80" = -5472
    "H82" = 2579.111
    "I82" = 404
    "J82" = 3666.6667
    "K82" = 1212
    "L82" = 11000.0001
    "M82" = -806
    "N82" = -11812.0001
    "H83" = 1200
    "I83" = 0
    "J83" = 1200
    "K83" = 0
    "L83" = 10800
    "N83" = -20160
    "H85" = 2579.111
    "I85" = 404
    "J85" = 3666.6667
    "K85" = 1212
    "L85" = 11000.0001
    "M85" = 192
    "N85" = -13808.0001
    "H94" = 2547.0588
    "I94" = 1466.6666
    "J94" = 2778.5715
    "K94" = 4399.9998
    "L94" = 8335.7145
    "M94" = -3723.9998
    "N94" = -9687.7145
    "H122" = 1150.3334
    "I122" = 504
    "J122" = 1279.6
    "K122" = 4536
    "L122" = 11516.4
    "M122" = -2086
    "N122" = -16416.4
    "H131" = 1461.7812
    "I131" = 451.25
    "J131" = 1606.1428
    "K131" = 1353.75
    "L131" = 4818.428400000001
    "M131" = 3686.25
    "N131" = -14898.4284
    "N76" = -11766.0001
    "N79" = -13652.0001
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$toClear = @("M80", "M83")
foreach ($ref in $toClear) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H122" = 927542.9399999999
    "I122" = 2778527.5
    "J122" = 2050.625
    "K122" = 8335582.5
    "L122" = 6151.875
    "M122" = -8333132.5
    "N122" = -11051.875
    "H126" = 2158.0881
    "I126" = 1404.2667
    "J126" = 2753.2104
    "K126" = 4212.800099999999
    "L126" = 8259.6312
    "M126" = -1742.800099999999
    "N126" = -13199.6312
    "H132" = 3474.7
    "I132" = 3227.375
    "J132" = 4464
    "K132" = 9682.125
    "L132" = 13392
    "M132" = -7152.125
    "N132" = -18452
    "H138" = 61680
    "J138" = 61680
    "L138" = 61680
    "N138" = -71960
    "H139" = 35884
    "J139" = 35884
    "L139" = 35884
    "N139" = -46164
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$values = @{
    "H16" = 16484.143
    "J16" = 50699
    "L16" = 50699
    "N16" = -51039
    "H132" = 4544.029
    "I132" = 3974.0588
    "J132" = 5082.3335
    "K132" = 11922.1764
    "L132" = 15247.0005
    "M132" = -9392.1764
    "N132" = -20307.0005
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$values = @{
    "H122" = 28368.895
    "I122" = 45117.87
    "K122" = 135353.61
    "M122" = -132903.61
    "H126" = 45245.176
    "I126" = 51831.95
    "J126" = 1333.3334
    "K126" = 155495.85
    "L126" = 4000.0002
    "M126" = -153025.85
    "N126" = -8940.0002
    "H131" = 67208.8
    "J131" = 67208.8
    "L131" = 67208.8
    "N131" = -77288.8
    "H132" = 16670992
    "I132" = 31255806
    "J132" = 2633.5715
    "K132" = 93767418
    "L132" = 7900.7145
    "M132" = -93764888
    "N132" = -12960.7145
    "H136" = 11942592
    "I136" = 19668138
    "J136" = 3112
    "K136" = 59004414
    "L136" = 9336
    "M136" = -59001864
    "N136" = -14436
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
